$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B column to "BBB" for the following rows
$rows_BBB = @(5, 13, 15, 19, 23, 29, 33, 34, 36, 45, 46, 47, 52, 55, 66, 71, 72, 80, 85, 86, 87, 96, 97, 98, 108, 113, 123, 125, 129, 131, 133, 135, 140, 142, 148, 152, 154, 158, 161, 165, 166, 167, 171, 172, 174, 175, 176, 180, 183, 184, 185, 189, 192, 195, 196, 199, 207, 209, 214, 219, 220, 221, 222, 224, 228, 229, 230, 238, 252, 256, 257, 258, 259, 260, 261, 262, 266, 267, 270, 271, 273, 274, 275, 279, 288, 295, 296, 301, 306, 307, 308, 309, 310, 314, 315, 318, 330, 332, 333, 337, 339, 340, 344, 346, 347, 352, 357, 361, 365, 366, 367, 370, 388, 389, 392, 397, 398, 399, 400, 403)
foreach ($r in $rows_BBB) {
    $ws.Cells.Item($r, 2).Value = "BBB"
}

# Set B column to "BB" for the following rows
$rows_BB = @(53, 74, 99, 100, 101, 109, 124, 126, 127, 130, 156, 157, 164, 182, 213, 241, 293, 294, 305, 323, 325, 327, 350, 378, 379, 380, 381, 382, 384, 396)
foreach ($r in $rows_BB) {
    $ws.Cells.Item($r, 2).Value = "BB"
}
